# Insert two new columns "m7" (D) and "m9" (E), pushing the existing
# repository/version/asc1/ac1/total columns two places to the right
# (D:H -> F:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift old columns D:H to F:J by inserting two blank columns before D.
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells, matching the style used by the other header cells
# (bold font, thin box border, centered / top aligned).
$ws.Range("D1").Value = "m7"
$ws.Range("E1").Value = "m9"
foreach ($addr in @("D1", "E1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

# Fill the m7 / m9 data columns for rows 2-15.
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 79.31
    $ws.Cells.Item($r, 5).Value = 8.050000000000001
}
